# Commit: "Fixed POI packaging and upgraded to POI 3.15."
#
# The underlying XML diff for this commit only reorders XML attributes
# and namespace (xmlns:*) declarations inside the start tags of
# word/document.xml, word/footer1.xml, word/footer2.xml, word/footer3.xml,
# word/footnotes.xml, word/header1.xml, word/header2.xml, word/header3.xml
# and word/styles.xml (e.g. <w:headerReference w:type="even" r:id="rId6"/>
# becoming <w:headerReference r:id="rId6" w:type="even"/>, <w:pgSz
# w:w="11906" w:h="16838"/> becoming <w:pgSz w:h="16838" w:w="11906"/>,
# etc.). Every attribute/value pair and every namespace URI present
# before the commit is still present after it, just alphabetised - this
# is exactly what happens when a document gets re-serialized by a newer
# OOXML writer (here: the POI 3.15 upgrade mentioned in the commit
# message) without anybody touching the document's actual content.
#
# No paragraph text, run/paragraph formatting, header/footer content,
# section properties values, relationships, or style definitions were
# added, removed, or changed by this commit - the rendered document and
# its Word object model are identical before and after. There is
# therefore no user-visible edit to apply through the Word object model
# here; this script intentionally performs no content mutation.
$d = $word.ActiveDocument
